# Apply cryptos list update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number;
# force text storage, then reset the style so no stray 's' attribute remains.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.862"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0900"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0997"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.52"
$ws.Range("D51").Style = "Normal"

# Plain text / percentage-string updates (safe to assign directly).
$ws.Range("D2").Value = "43.321.92"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "2.282.20"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +4.02%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "2.625.33"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "2.281.88"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").Value = "43.352.53"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("E35").Value = "  +3.68%  "
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  +8.69%  "
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("E51").Value = "  +28.14%  "
